$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bsky accounts tracked -> four new columns appended after Z (AA:AD)
$ws.Range("AA1").Value = "fuelpovertyaction.bsky.social"
$ws.Range("AB1").Value = "jrf-uk.bsky.social"
$ws.Range("AC1").Value = "e3g.bsky.social"
$ws.Range("AD1").Value = "neweconomics.bsky.social"

# Corrected/updated post counts that were previously recorded as 0
$ws.Range("N12").Value = 97
$ws.Range("P12").Value = 750
$ws.Range("H13").Value = 401

# Match the author's final selection (new columns highlighted) before saving
$ws.Range("AA1:AD1").Select()
